# Generate Report for Handback
# Updates the localization-status workbook to reflect a completed handback:
#  - Status text changes from "In Translation" to "Handed back: in sync with en-US"
#    on the Overview sheet (zh-cn/de-de columns) and on each locale sheet's Status column.
#  - Each locale sheet's "Latest Target File" / "Latest Handback File" /
#    "Latest Handback DateTime" columns get populated for both rows.
#  - A hyperlink (pointing at the same source .md on GitHub as column A) is added
#    on the new "Latest Target File" cells.
#  - A few columns are widened so the new, longer text fits.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$mdTarget1 = "6af204bc-fcdf-4a43-a2c7-644927def291.md"
$mdTarget2 = "a54c2a7b-2b00-4316-91f6-aa67d3d260ac.md"
$mdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d536d7904feab985507d6dd0a73a897405d85533/e2e/6af204bc-fcdf-4a43-a2c7-644927def291.md"
$mdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d536d7904feab985507d6dd0a73a897405d85533/e2e/a54c2a7b-2b00-4316-91f6-aa67d3d260ac.md"

# ---------------------------------------------------------------------------
# Overview sheet: the zh-cn / de-de status columns (E, F) show the same text
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsOverview.Columns.Item(5).ColumnWidth = 29.9777047293527
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777047293527

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsZh.Range("I2").Value = $mdTarget1
$wsZh.Range("J2").Value = "6af204bc-fcdf-4a43-a2c7-644927def291.1b83f624619bac49c7271d419b703a352ad4d964.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-09-06 16:48:00"

$wsZh.Range("I3").Value = $mdTarget2
$wsZh.Range("J3").Value = "a54c2a7b-2b00-4316-91f6-aa67d3d260ac.0b199ada4f0026388de36fbae71bf61a6369bff6.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-09-06 16:48:00"

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $mdUrl1, "", "", $mdTarget1) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $mdUrl2, "", "", $mdTarget2) | Out-Null

$wsZh.Range("I2").Font.Underline = $true
$wsZh.Range("I2").Font.Color = $wsZh.Range("A2").Font.Color
$wsZh.Range("I3").Font.Underline = $true
$wsZh.Range("I3").Font.Color = $wsZh.Range("A3").Font.Color

$wsZh.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsZh.Columns.Item(9).ColumnWidth = 40
$wsZh.Columns.Item(10).ColumnWidth = 40

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

$wsDe.Range("I2").Value = $mdTarget1
$wsDe.Range("J2").Value = "6af204bc-fcdf-4a43-a2c7-644927def291.1b83f624619bac49c7271d419b703a352ad4d964.de-de.xlf"
$wsDe.Range("K2").Value = "2016-09-06 16:48:21"

$wsDe.Range("I3").Value = $mdTarget2
$wsDe.Range("J3").Value = "a54c2a7b-2b00-4316-91f6-aa67d3d260ac.0b199ada4f0026388de36fbae71bf61a6369bff6.de-de.xlf"
$wsDe.Range("K3").Value = "2016-09-06 16:48:21"

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $mdUrl1, "", "", $mdTarget1) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $mdUrl2, "", "", $mdTarget2) | Out-Null

$wsDe.Range("I2").Font.Underline = $true
$wsDe.Range("I2").Font.Color = $wsDe.Range("A2").Font.Color
$wsDe.Range("I3").Font.Underline = $true
$wsDe.Range("I3").Font.Color = $wsDe.Range("A3").Font.Color

$wsDe.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsDe.Columns.Item(9).ColumnWidth = 40
$wsDe.Columns.Item(10).ColumnWidth = 40

Write-Host "Handback report generated."
